# Update "想去人数" (want-to-go count) values in column F
# for sheets "展览" (Exhibition) and "全部类型" (All types).
# Rows 4, 7, 13, 15 are intentionally left unchanged.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1185
    3  = 601
    5  = 35
    6  = 180
    8  = 64
    9  = 10
    10 = 5528
    11 = 4917
    12 = 20
    14 = 5
    16 = 202
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
